# Update coin price (D) / volume-change (E) figures with the latest refresh values.
# Cells in this sheet are stored as literal text (e.g. "314.66", "2.87%"), not numbers,
# so each write forces a text NumberFormat, then restores the default "Normal" style
# to avoid leaving a stray per-cell format behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "314.66"
    "E2" = "2.87%"
    "E3" = "-2.29%"
    "D4" = "5.133"
    "E4" = "0.60%"
    "E5" = "3.07%"
    "D6" = "2.124"
    "E6" = "-0.40%"
    "D7" = "7.994"
    "E7" = "0.45%"
    "D8" = "4.151"
    "E8" = "0.89%"
    "D9" = "0.9286"
    "E9" = "1.24%"
    "D10" = "0.1002"
    "E10" = "3.80%"
    "D11" = "0.1871"
    "D12" = "0.09103"
    "E12" = "4.92%"
    "D13" = "0.03604"
    "E13" = "1.53%"
    "D14" = "0.09908"
    "E14" = "-0.29%"
    "D15" = "0.001430"
    "E15" = "-0.59%"
    "D16" = "0.005713"
    "E16" = "-0.35%"
    "D17" = "3.467"
    "E17" = "-0.08%"
    "D18" = "2.755"
    "E18" = "1.55%"
    "D19" = "0.3410"
    "E19" = "0.49%"
    "D20" = "0.1331"
    "D21" = "5.096"
    "E21" = "-1.37%"
    "E22" = "9.85%"
    "D23" = "0.04557"
    "E23" = "-0.06%"
    "D25" = "0.004704"
    "E25" = "-6.66%"
    "E26" = "-21.92%"
    "D27" = "0.0004508"
    "E27" = "-5.14%"
    "D39" = "0.01946"
    "E39" = "4.96%"
    "D40" = "0.04850"
    "E40" = "1.67%"
    "D41" = "0.007725"
    "E41" = "2.77%"
    "D42" = "0.1392"
    "E42" = "-0.61%"
    "D43" = "0.007851"
    "E43" = "1.37%"
    "D44" = "0.002112"
    "E44" = "-5.34%"
    "D45" = "0.01179"
    "E45" = "6.72%"
    "D46" = "0.00006612"
    "E46" = "4.36%"
    "E47" = "0.05%"
    "D48" = "39.28"
    "E48" = "-17.23%"
    "D49" = "0.001703"
    "E49" = "-14.90%"
    "E50" = "0.05%"
    "E51" = "0.05%"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    $cell.Style = "Normal"
}
